$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rename the first CTE from "all_journeys" to "filtered_dp013" on the
#    "WITH ..." line.
# ---------------------------------------------------------------------------
$withPara = $d.Paragraphs.Item(6)
$withPara.Range.Text = "WITH filtered_dp013 AS ("

# ---------------------------------------------------------------------------
# 2. Insert the body of the (soon to be relocated) dp013 subquery right after
#    the "WITH filtered_dp013 AS (" line, followed by the closing "),"
#    and the original "all_journeys AS (" line.
# ---------------------------------------------------------------------------
$withPara.Range.InsertParagraphAfter()
$newBodyPara = $d.Paragraphs.Item(7)
$newLines = @(
    "    SELECT traceability.*, ",
    "           ROW_NUMBER() OVER (PARTITION BY cdm_join_key ORDER BY request_time DESC) AS last_record",
    "    FROM ``bt-bvp-ml-plat-ai-pipe-exp.hypothesis_testing.dp013_acquisition`` AS traceability",
    "    WHERE LOWER(traceability.product_type) = 'broadband'",
    "      AND DATE(request_time) = DATE '2024-12-20'",
    "      AND cdm_join_key IS NOT NULL",
    "    QUALIFY last_record = 1",
    "),",
    "all_journeys AS ("
)
$newBodyPara.Range.Text = [string]::Join([string][char]13, $newLines)

# ---------------------------------------------------------------------------
# 3. Replace the old inline subquery ("FROM ( ... ) QUALIFY ... AS dp013")
#    with a simple reference to the new "filtered_dp013" CTE.
# ---------------------------------------------------------------------------
$fromParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "FROM (`r") {
        $fromParaIndex = $i
        break
    }
}

$qualifyParaIndex = $null
for ($i = $fromParaIndex; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "QUALIFY last_record = 1 AS dp013`r") {
        $qualifyParaIndex = $i
        break
    }
}

$fromPara = $d.Paragraphs.Item($fromParaIndex)
$qualifyPara = $d.Paragraphs.Item($qualifyParaIndex)

$oldSubquery = $d.Range($fromPara.Range.Start, $qualifyPara.Range.End)
$oldSubquery.Delete()

$d.Paragraphs.Item($fromParaIndex).Range.Text = "FROM filtered_dp013 AS dp013"
